$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.053.79'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '3.202.06'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.29'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.201.20'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.515'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.68%  '
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '39.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.59%  '
$ws.Range('D15').Value = '3.726.00'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '66.014.65'
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('D18').Value = '3.201.58'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '512.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.742'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  +5.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('E29').Value = '  +3.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.91'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.77%  '
$ws.Range('E31').Value = '  +3.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.57'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.04'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '485.42'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.15%  '
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.92'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('E42').Value = '  +3.42%  '
$ws.Range('E43').Value = '  +4.91%  '
$ws.Range('D44').Value = '0.0₃0651'
$ws.Range('E44').Value = '  +9.45%  '
$ws.Range('D45').Value = '2.942.18'
$ws.Range('E45').Value = '  -3.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.43'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.86%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  +1.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.61%  '
